# 027 Week 15/16 update
# Fill in WK 16 (column S) scores for several pairs on Sheet1.
# Downstream totals (AB/AC on Sheet1) and the COUNTIF "games played"
# helpers on the hidden "xxDO NOT EDITxx" sheet recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$scores = @{
    "S12" = 31
    "S13" = 35
    "S16" = 31
    "S17" = 23
    "S18" = 28
    "S20" = 28
    "S26" = 30
    "S28" = 31
    "S29" = 38
    "S31" = 30
}

foreach ($addr in $scores.Keys) {
    $ws.Range($addr).Value = $scores[$addr]
}
